$wb = $excel.ActiveWorkbook

# ==== ALC ====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 5702.3335
$ws.Range("I82").Value = 1488.2
$ws.Range("J82").Value = 8712.429
$ws.Range("K82").Value = 4464.6
$ws.Range("L82").Value = 26137.287
$ws.Range("M82").Value = -4058.6
$ws.Range("N82").Value = -26949.287
$ws.Range("H85").Value = 5702.3335
$ws.Range("I85").Value = 1488.2
$ws.Range("J85").Value = 8712.429
$ws.Range("K85").Value = 4464.6
$ws.Range("L85").Value = 26137.287
$ws.Range("M85").Value = -3060.6
$ws.Range("N85").Value = -28945.287
$ws.Range("H112").Value = 1322.6852
$ws.Range("J112").Value = 1322.6852
$ws.Range("L112").Value = 3968.0556
$ws.Range("N112").Value = -6184.0556
$ws.Range("H118").Value = 1010
$ws.Range("J118").Value = 2000
$ws.Range("L118").Value = 6000
$ws.Range("N118").Value = -9314
$ws.Range("H132").Value = 34488324
$ws.Range("I132").Value = 38466900
$ws.Range("K132").Value = 115400700
$ws.Range("M132").Value = -115398170

# ==== ARM ====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 944.375
$ws.Range("I2").Value = 755.5
$ws.Range("J2").Value = 1133.25
$ws.Range("K2").Value = 755.5
$ws.Range("L2").Value = 1133.25
$ws.Range("M2").Value = -642.5
$ws.Range("N2").Value = -1359.25
$ws.Range("H36").Value = 7074.5
$ws.Range("I36").Value = 3113
$ws.Range("K36").Value = 3113
$ws.Range("M36").Value = -2767
$ws.Range("H45").Value = 5500
$ws.Range("I45").Value = 5500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5500
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -5123
$ws.Range("H61").Value = 1733
$ws.Range("I61").Value = 1688.2
$ws.Range("J61").Value = 1845
$ws.Range("K61").Value = 1688.2
$ws.Range("L61").Value = 1845
$ws.Range("M61").Value = -1476.2
$ws.Range("N61").Value = -2269
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -814
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 6000
$ws.Range("M89").Value = -72
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H116").Value = 944.375
$ws.Range("I116").Value = 755.5
$ws.Range("J116").Value = 1133.25
$ws.Range("K116").Value = 755.5
$ws.Range("L116").Value = 1133.25
$ws.Range("M116").Value = 1538.5
$ws.Range("N116").Value = -5721.25
$ws.Range("H133").Value = 31260
$ws.Range("J133").Value = 31260
$ws.Range("L133").Value = 31260
$ws.Range("N133").Value = -36320
$ws.Range("H136").Value = 1733
$ws.Range("I136").Value = 1688.2
$ws.Range("J136").Value = 1845
$ws.Range("K136").Value = 5064.6
$ws.Range("L136").Value = 5535
$ws.Range("M136").Value = -2514.6
$ws.Range("N136").Value = -10635
$ws.Range("H139").Value = 42946.562
$ws.Range("J139").Value = 42946.562
$ws.Range("L139").Value = 42946.562
$ws.Range("N139").Value = -53226.562

# ==== BSM ====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 944.375
$ws.Range("I3").Value = 755.5
$ws.Range("J3").Value = 1133.25
$ws.Range("K3").Value = 755.5
$ws.Range("L3").Value = 1133.25
$ws.Range("M3").Value = -641.5
$ws.Range("N3").Value = -1361.25
$ws.Range("H7").Value = 3096055
$ws.Range("J7").Value = 3656475
$ws.Range("L7").Value = 3656475
$ws.Range("N7").Value = -3656701
$ws.Range("H9").Value = 21850
$ws.Range("J9").Value = 21850
$ws.Range("L9").Value = 21850
$ws.Range("N9").Value = -22186
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws.Range("H107").Value = 2069.625
$ws.Range("I107").Value = 1900.1111
$ws.Range("J107").Value = 2287.5715
$ws.Range("K107").Value = 1900.1111
$ws.Range("L107").Value = 2287.5715
$ws.Range("M107").Value = 19.88889999999992
$ws.Range("N107").Value = -6127.5715
$ws.Range("H134").Value = 3077.724
$ws.Range("I134").Value = 1658.6666
$ws.Range("K134").Value = 4975.9998
$ws.Range("M134").Value = -2440.9998

# ==== CRP ====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 6448.353
$ws.Range("J12").Value = 6712.3125
$ws.Range("L12").Value = 6712.3125
$ws.Range("N12").Value = -7052.3125
$ws.Range("H58").Value = 2069.4517
$ws.Range("J58").Value = 4151.3335
$ws.Range("L58").Value = 4151.3335
$ws.Range("N58").Value = -4557.3335
$ws.Range("H82").Value = 39800
$ws.Range("J82").Value = 39800
$ws.Range("L82").Value = 39800
$ws.Range("N82").Value = -40522
$ws.Range("H85").Value = 39800
$ws.Range("J85").Value = 39800
$ws.Range("L85").Value = 39800
$ws.Range("N85").Value = -42296
$ws.Range("H112").Value = 34846.152
$ws.Range("J112").Value = 34846.152
$ws.Range("L112").Value = 34846.152
$ws.Range("N112").Value = -37800.152
$ws.Range("H119").Value = 40761
$ws.Range("J119").Value = 40761
$ws.Range("L119").Value = 40761
$ws.Range("N119").Value = -50437
$ws.Range("H136").Value = 2069.4517
$ws.Range("J136").Value = 4151.3335
$ws.Range("L136").Value = 12454.0005
$ws.Range("N136").Value = -17554.0005

# ==== CUL ====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 743373.2
$ws.Range("I5").Value = 850
$ws.Range("J5").Value = 1028959
$ws.Range("K5").Value = 2550
$ws.Range("L5").Value = 3086877
$ws.Range("M5").Value = -2438
$ws.Range("N5").Value = -3087101
$ws.Range("H96").Value = 12759
$ws.Range("J96").Value = 12759
$ws.Range("L96").Value = 38277
$ws.Range("N96").Value = -42395
$ws.Range("H135").Value = 743373.2
$ws.Range("I135").Value = 850
$ws.Range("J135").Value = 1028959
$ws.Range("K135").Value = 7650
$ws.Range("L135").Value = 9260631
$ws.Range("M135").Value = -5115
$ws.Range("N135").Value = -9265701

# ==== GSM ====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5046
$ws.Range("I132").Value = 3351.7144
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 10055.1432
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -7525.143199999999
$ws.Range("N132").Value = -32057.999

# ==== LTW ====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8752.666999999999
$ws.Range("I132").Value = 3656.25
$ws.Range("J132").Value = 12829.8
$ws.Range("K132").Value = 10968.75
$ws.Range("L132").Value = 38489.39999999999
$ws.Range("M132").Value = -8438.75
$ws.Range("N132").Value = -43549.39999999999
$ws.Range("H136").Value = 3409.6562
$ws.Range("I136").Value = 1208.3684
$ws.Range("J136").Value = 6626.923
$ws.Range("K136").Value = 3625.1052
$ws.Range("L136").Value = 19880.769
$ws.Range("M136").Value = -1075.1052
$ws.Range("N136").Value = -24980.769

# ==== WVR ====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 34400
$ws.Range("J82").Value = 34400
$ws.Range("L82").Value = 34400
$ws.Range("N82").Value = -35166
$ws.Range("H85").Value = 34400
$ws.Range("J85").Value = 34400
$ws.Range("L85").Value = 34400
$ws.Range("N85").Value = -37052
$ws.Range("H112").Value = 35500
$ws.Range("J112").Value = 35500
$ws.Range("L112").Value = 35500
$ws.Range("N112").Value = -38454
$ws.Range("H113").Value = 7432.357
$ws.Range("I113").Value = 9370
$ws.Range("J113").Value = 327.66666
$ws.Range("K113").Value = 28110
$ws.Range("L113").Value = 982.9999799999999
$ws.Range("M113").Value = -25940
$ws.Range("N113").Value = -5322.99998
$ws.Range("H115").Value = 39192.31
$ws.Range("J115").Value = 39192.31
$ws.Range("L115").Value = 39192.31
$ws.Range("N115").Value = -42326.31
$ws.Range("H118").Value = 29890
$ws.Range("J118").Value = 29890
$ws.Range("L118").Value = 29890
$ws.Range("N118").Value = -33204
$ws.Range("H136").Value = 4870.522
$ws.Range("I136").Value = 4592.8965
$ws.Range("J136").Value = 5344.1177
$ws.Range("K136").Value = 13778.6895
$ws.Range("L136").Value = 16032.3531
$ws.Range("M136").Value = -11228.6895
$ws.Range("N136").Value = -21132.3531
